$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 4) duplicating the Daphne/Boone contact, with a new
# inbox webhook email address.
$ws.Range("C4").Value = "Daphne"
$ws.Range("G4").Value = "Boone"
$ws.Range("H4").Value = "shahimran@outlook.com"

# Turn H4 into a mailto: hyperlink (adds the relationship + hyperlink entry).
$ws.Hyperlinks.Add($ws.Range("H4"), "mailto:shahimran@outlook.com")

# Hyperlinks.Add re-applies hyperlink formatting with a fresh style index;
# put it back in line with the existing hyperlink cells (H2/H3) so H4 shares
# the same style as the rest of the EMAIL column.
$ws.Range("H4").Style = $ws.Range("H3").Style

# Move the active selection to G7, as recorded in the saved sheet view.
$ws.Range("G7").Select()
